$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 162, shifting existing rows 162:237 down to 163:238
$ws.Rows("162:162").Insert()

# Populate the newly inserted row 162 with its data
$ws.Cells.Item(162, 1).Value2 = 9
$ws.Cells.Item(162, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(162, 3).Value2 = "Metropolitana"
$ws.Cells.Item(162, 4).Value2 = 44466
$ws.Cells.Item(162, 5).Value2 = 13
$ws.Cells.Item(162, 6).Value2 = 100112012
$ws.Cells.Item(162, 7).Value2 = "Espinaca"
$ws.Cells.Item(162, 8).Value2 = "Sin especificar"
$ws.Cells.Item(162, 9).Value2 = "Primera"
$ws.Cells.Item(162, 10).Value2 = 131
$ws.Cells.Item(162, 11).Value2 = 7000
$ws.Cells.Item(162, 12).Value2 = 8000
$ws.Cells.Item(162, 13).Value2 = 7496
$ws.Cells.Item(162, 14).Value2 = "$/cuna 10 kilos"
$ws.Cells.Item(162, 15).Value2 = "Provincia de Chacabuco"
$ws.Cells.Item(162, 16).Value2 = 750
$ws.Cells.Item(162, 17).Value2 = 10
$ws.Cells.Item(162, 18).Value2 = "Hortaliza"
